# Add a new "Title Only" slide at the end of the deck (PpSlideLayout 11 = ppLayoutTitleOnly),
# matching slideLayout6.xml ("Titre seul"), then set its title text.
$p = $ppt.ActivePresentation

$newSlide = $p.Slides.Add($p.Slides.Count + 1, 11)

$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Blank slide to be ignored"
